$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right
#    after the title (Heading1) paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Re-insert a bold "Play Celebrity Zoo Free - Unique Animal Themed
#    Slot Game" paragraph right before the final (DALLE prompt) paragraph.
#    We do this by replacing the last "real" bullet paragraph
#    ("Limited variety in bonus games") with itself followed by the new
#    bold paragraph, using InsertXML so the run structure comes out
#    exactly as Word would author it (leading empty run + bold run).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastBullet = $d.Paragraphs.Item($count - 1)

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$bulletXml = '<w:p ' + $xmlNs + '><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Limited variety in bonus games</w:t></w:r></w:p>'
$newParaXml = '<w:p ' + $xmlNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Celebrity Zoo Free - Unique Animal Themed Slot Game</w:t></w:r></w:p>'

$r = $lastBullet.Range.Duplicate
$r.Collapse(1)
$r.InsertXML($bulletXml + $newParaXml)

# ------------------------------------------------------------------
# 3) Turn the old DALLE-prompt paragraph text into the new meta
#    description text (it keeps its italic run formatting).
# ------------------------------------------------------------------
$sq = [char]39
$oldPrompt = 'Prompt: "Create a cartoon-style image for the game ' + $sq + 'Celebrity Zoo' + $sq + ' that features a happy Maya warrior with glasses." DALLE, I need you to create a feature image for the slot game ' + $sq + 'Celebrity Zoo' + $sq + ' that highlights its quirky and playful atmosphere. The image should be in a cartoon style that features a happy Maya warrior with glasses. The Maya warrior should be depicted in bright colors to match the game' + $sq + 's whimsical graphics, with a big smile on his face to convey the game' + $sq + 's fun and lighthearted tone. The Maya warrior should also be holding a camera to symbolize the game' + $sq + 's Paparazzi Bonus Game. This feature image should be eye-catching and vibrant to catch the players' + $sq + ' attention and make them want to try out the game. Make sure the image is designed to fit the game' + $sq + 's theme and can convey the game' + $sq + 's exciting features.'
$newDescription = "Enjoy the fun and engaging gameplay of Celebrity Zoo, an online slot game featuring unique animal characters and exciting bonus games. Play now for free!"

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)
